$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("survey")

# ---------------------------------------------------------------------------
# The "label features" demo block (rows 14-17) is being reorganised:
#   - the HTML-label / HTML-hint note moves up to row 14
#   - the HTML "img/audio" note (with its audio+image columns) moves to row 15
#   - the "Enter your name" text question moves to row 16, with a new hint
#   - the Handlebars-template note moves to row 17, with an expanded label
# ---------------------------------------------------------------------------

# --- fix up formatting (wrap text) on the cells that need it before the
#     values get shuffled around, copying from cells that already carry the
#     target formatting ---
$ws.Range("D16").Copy() | Out-Null
$ws.Range("D14").PasteSpecial(-4122) | Out-Null   # xlPasteFormats -> no-wrap style on D14
$ws.Range("E16").Copy() | Out-Null
$ws.Range("E14").PasteSpecial(-4122) | Out-Null   # xlPasteFormats -> no-wrap style on E14

$ws.Range("G17").Copy() | Out-Null
$ws.Range("G15").PasteSpecial(-4122) | Out-Null   # xlPasteFormats -> matches audio column style
$ws.Range("H17").Copy() | Out-Null
$ws.Range("H15").PasteSpecial(-4122) | Out-Null   # xlPasteFormats -> matches image column style

$ws.Range("C14").Copy() | Out-Null
$ws.Range("C16").PasteSpecial(-4122) | Out-Null   # xlPasteFormats -> matches name column style

# --- now move the actual content ---
$ws.Range("A14").Value = "note"
$ws.Range("C14").ClearContents()
$ws.Range("D14").Value = '<u>Labels</u> <i>can</i> contain <span style="color:red;">HTML</span>'
$ws.Range("E14").Value = "So can <b>hints</b>"

$ws.Range("A15").Value = "note"
$ws.Range("D15").Value = "Labels can contain &lt;img&gt; and &lt;audio&gt; HTML tags, but it is often easier to add media via the image and audio columns."
$ws.Range("E15").ClearContents()
$ws.Range("G15").Value = "audio/carrioncrow.mp3"
$ws.Range("H15").Value = "img/dolphin.png"

$ws.Range("A16").Value = "text"
$ws.Range("C16").Value = "name"
$ws.Range("D16").Value = "Enter your name"
$ws.Range("E16").Value = "It will be used in the next question."

$ws.Range("A17").Value = "note"
$ws.Range("D17").Value = "<h3>`nThis label uses Handlesbars template features:`n</h3>`n{{#if name}}`nHello {{name}}!`n{{else}}`nName not entered.`n{{/if}}"
$ws.Range("E17").Value = "Handlebars templates allow labels to change depending on the values previously entered."
$ws.Range("G17").ClearContents()
$ws.Range("H17").ClearContents()

# ---------------------------------------------------------------------------
# The "born" date question's Happy-Birthday note: guard the condition against
# an empty/undefined 'born' value, and add an explanatory hint.
# ---------------------------------------------------------------------------
$ws.Range("B23").Value = "data('born') && data('born').getDay() === now().getDay() && data('born').getMonth() === now().getMonth()"
$ws.Range("E23").Value = "This prompt shows how to use dates in fomulas."

# The "time" question's label is reworded.
$ws.Range("D24").Value = "What time do you usually wake up?"

# ---------------------------------------------------------------------------
# Settings sheet: drop the custom-appearance settings (font-size / theme) -
# that demo is being moved out into its own, separate form.
# ---------------------------------------------------------------------------
$settings = $wb.Worksheets.Item("settings")
$settings.Range("A5:B6").EntireRow.Delete() | Out-Null
